$wb = $excel.ActiveWorkbook
$wsCases = $wb.Worksheets.Item("CasesUSA")
$wsFatal = $wb.Worksheets.Item("FatalitiesUSA")

# --- CasesUSA (sheet1) updates ---

# Updated parameter inputs in column D (rows 1-3)
$wsCases.Range("D1").Value = 0.0055
$wsCases.Range("D2").Value = 0.3974
$wsCases.Range("D3").Value = 3.4251

# Corrected case count for row 33
$wsCases.Range("A33").Value = 25715

# Append a new data point (row 34), copying the formatting from row 33
# so the new cells keep the same style as the rest of the table.
$wsCases.Range("A33:B33").Copy()
$wsCases.Range("A34").PasteSpecial(-4122)
$wsCases.Range("A34").Value = 30459
$wsCases.Range("B34").Formula = "=LN(A34)"

# --- View / selection state ---
# Scroll/select on FatalitiesUSA first, then finish on CasesUSA so that
# CasesUSA ends up as the active (tabSelected) sheet, matching the saved
# workbook state.
$wsFatal.Range("A11").Select()
$wsFatal.Range("A30").Select()

$wsCases.Activate()
$wsCases.Range("D4").Select()
